$d = $word.ActiveDocument

# 1) Brief paragraph
$d.Content.Find.Execute(
    "An email sent to partners who have attended the event. This email will include a photo gallery It will be sent via customer.io",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Một email gửi đến các đối tác đã tham dự sự kiện. Email này sẽ bao gồm một thư viện ảnh. Nó sẽ được gửi qua customer.io.",
    2) | Out-Null

# 2) Target audience paragraph
$d.Content.Find.Execute(
    "Event attendees",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Người tham dự sự kiện",
    2) | Out-Null

# 3) Subject line lead-in text
$d.Content.Find.Execute(
    "Thank you for coming to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cảm ơn bạn đã tham dự sự kiện ",
    2) | Out-Null

# 3b) [EVENT NAME] inside the Subject line only (highlighted run) -
#     the later "[EVENT NAME]" placeholder in the body must stay untouched,
#     so scope the Find to the Subject paragraph only.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Subject:*") {
        $p.Range.Find.Execute(
            "[EVENT NAME]",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "[TÊN SỰ KIỆN]",
            2) | Out-Null
        break
    }
}

# 4) "You made our event a success! ..." heading
$d.Content.Find.Execute(
    "You made our event a success! 🎉",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sự tham dự của bạn đã góp phần vào thành công của sự kiện chúng tôi! 🎉",
    2) | Out-Null

# 5) Greeting
$d.Content.Find.Execute(
    "Hi ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Xin chào ",
    2) | Out-Null

# 6) "Thank you for attending " lead-in
$d.Content.Find.Execute(
    "Thank you for attending ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cảm ơn bạn đã tham dự sự kiện ",
    2) | Out-Null

# 7) " in " connector between [EVENT NAME] and [CITY]
$d.Content.Find.Execute(
    " in ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " tại ",
    2) | Out-Null

# 8) Closing sentence after [COUNTRY]
$d.Content.Find.Execute(
    ". We hope you had a great time, and it was a pleasure getting to know you!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Chúng tôi hy vọng bạn đã có một khoảng thời gian tuyệt vời. Rất vinh dự khi được làm quen với bạn!",
    2) | Out-Null

# 9) Closing paragraph near the end of the document
$d.Content.Find.Execute(
    "We hope the event inspired you as much as it did us, and let’s keep growing together!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Chúng tôi mong rằng sự kiện đã có thể truyền cảm hứng cho bạn như nó đã làm với chúng tôi. Chúng tôi hy vọng bạn và chúng tôi sẽ cùng nhau phát triển hơn nữa trong tương lai!",
    2) | Out-Null
